$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# dataprovider refresh: replace the "Show Number" contact values in
# column B (rows 2-6) with freshly generated phone numbers.
# Force text formatting first so the leading zeros in the phone
# numbers are preserved (otherwise Excel would coerce the numeric-
# looking strings into numbers), then clear the temporary format so
# the cells fall back to the sheet's default (General) style.
$rng = $ws.Range("B2:B6")
$rng.NumberFormat = "@"

$ws.Range("B2").Value = "07947122898"
$ws.Range("B3").Value = "07947108658"
$ws.Range("B4").Value = "07947119122"
$ws.Range("B5").Value = "07947119177"
$ws.Range("B6").Value = "07942700016"

$rng.ClearFormats()
